# Updated capital structure database — refresh computed financial metrics
# for the three Taiwan regional-bank rows (2-4) on the single worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (aggregate / row "2") ---
$ws.Range("D2").Value = 0.0505
$ws.Range("E2").Value = 0.0345
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 33.99
$ws.Range("L2").Value = 0.2170498084291188
$ws.Range("M2").Value = 6.140000000000001
$ws.Range("N2").Value = 0.01191076624636276
$ws.Range("O2").Value = 0.1806413651073845
$ws.Range("P2").Value = 6.140000000000001
$ws.Range("Q2").Value = 0.01191076624636276
$ws.Range("R2").Value = 0.1806413651073845
$ws.Range("U2").Value = 279.3
$ws.Range("V2").Value = 0.54180407371484
$ws.Range("W2").Value = 0.04733683338802169
$ws.Range("X2").Value = 0.1179053464609755
$ws.Range("Y2").Value = -0.07056851307295378
$ws.Range("Z2").Value = 0.09407665505226481
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04149570876818139
$ws.Range("AC2").Value = -0.04149570876818139
$ws.Range("AD2").Value = 1405.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1405.7
$ws.Range("AG2").Value = 1126.4
$ws.Range("AH2").Value = 0.7316781178430147
$ws.Range("AI2").Value = 0.6548800372699743
$ws.Range("AJ2").Value = 0.6860344722577502
$ws.Range("AK2").Value = 0.6032562125107113

# --- Row 3 (Bank of Kaohsiung Co., Ltd.) ---
$ws.Range("D3").Value = 0.0421
$ws.Range("E3").Value = 0.0449
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 27
$ws.Range("L3").Value = 0.2290076335877863
$ws.Range("M3").Value = 3.95
$ws.Range("N3").Value = 0.009823427008206913
$ws.Range("O3").Value = 0.1462962962962963
$ws.Range("P3").Value = 3.95
$ws.Range("Q3").Value = 0.009823427008206913
$ws.Range("R3").Value = 0.1462962962962963
$ws.Range("U3").Value = 222.2
$ws.Range("V3").Value = 0.5525988560059686
$ws.Range("W3").Value = 0.05484460694698354
$ws.Range("X3").Value = 0.0872168410126459
$ws.Range("Y3").Value = -0.03237223406566236
$ws.Range("Z3").Value = 0.09702905110690478
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04116001090918232
$ws.Range("AC3").Value = -0.04116001090918232
$ws.Range("AD3").Value = 860.3
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 860.3
$ws.Range("AG3").Value = 638.0999999999999
$ws.Range("AH3").Value = 0.6814797211660328
$ws.Range("AI3").Value = 0.6077711056163899
$ws.Range("AJ3").Value = 0.6134397231301673
$ws.Range("AK3").Value = 0.5347356071398642

# --- Row 4 (Taipei Star Bank) ---
$ws.Range("D4").Value = 0.05889999999999999
$ws.Range("E4").Value = 0.0241
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6.99
$ws.Range("L4").Value = 0.1806201550387597
$ws.Range("M4").Value = 2.19
$ws.Range("N4").Value = 0.01931216931216931
$ws.Range("O4").Value = 0.3133047210300429
$ws.Range("P4").Value = 2.19
$ws.Range("Q4").Value = 0.01931216931216931
$ws.Range("R4").Value = 0.3133047210300429
$ws.Range("U4").Value = 57.1
$ws.Range("V4").Value = 0.5035273368606702
$ws.Range("W4").Value = 0.03982905982905983
$ws.Range("X4").Value = 0.148593851909305
$ws.Range("Y4").Value = -0.1087647920802452
$ws.Range("Z4").Value = 0.08609566184649611
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04183140662718045
$ws.Range("AC4").Value = -0.04183140662718045
$ws.Range("AD4").Value = 545.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 545.4
$ws.Range("AG4").Value = 488.3
$ws.Range("AH4").Value = 0.8278688524590164
$ws.Range("AI4").Value = 0.7461012311901505
$ws.Range("AJ4").Value = 0.8115339870367293
$ws.Range("AK4").Value = 0.7245882178364742

# debt_ebitda / net_debt_ebitda (columns AN, AP) are no longer populated
# for these rows in the refreshed database — remove the stale values.
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
